# Updated capital structure database
# czech_republic_bank_money_center.xlsx - refresh the three Czech "Bank (Money
# Center)" rows (anonymized company, MONETA Money Bank, Komercni banka) with
# new capital-structure figures. Komercni banka and MONETA also swap row
# order (row 3 <-> row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (anonymized "2" company) ---------------------------------------
$ws.Range("D2").Value  = -0.05455
$ws.Range("E2").Value  = -0.07325000000000001
$ws.Range("F2").Value  = -0.0243
$ws.Range("I2").Value  = 0
$ws.Range("J2").Value  = 0
$ws.Range("K2").Value  = 555.3
$ws.Range("L2").Value  = 0.3660514172709294
$ws.Range("M2").Value  = 73.2
$ws.Range("N2").Value  = 0.009894031141868513
$ws.Range("O2").Value  = 0.1318206374932469
$ws.Range("P2").Value  = 73.2
$ws.Range("Q2").Value  = 0.009894031141868513
$ws.Range("R2").Value  = 0.1318206374932469
$ws.Range("U2").Value  = 2532.5
$ws.Range("V2").Value  = 0.342303741349481
$ws.Range("W2").Value  = 0.108641618237655
$ws.Range("X2").Value  = 0.05853195823809115
$ws.Range("Y2").Value  = 0.05010965999956389
$ws.Range("Z2").Value  = 0.1029402783526841
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.03960677643099839
$ws.Range("AC2").Value = -0.03960677643099839
$ws.Range("AD2").Value = 9589.799999999999
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 9589.799999999999
$ws.Range("AG2").Value = 7057.299999999999
$ws.Range("AH2").Value = 0.5644977101752982
$ws.Range("AI2").Value = 0.6102050815426611
$ws.Range("AJ2").Value = 0.4882018857613259
$ws.Range("AK2").Value = 0.535325262455246
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# --- Row 3: now Komercni banka, a.s. (SEP:KOMB) (was MONETA) --------------
$ws.Range("B3").Value  = "Komercní banka, a.s. (SEP:KOMB)"
$ws.Range("D3").Value  = -0.016
$ws.Range("E3").Value  = -0.0557
$ws.Range("F3").Value  = -0.0243
$ws.Range("I3").Value  = 0
$ws.Range("J3").Value  = 0
$ws.Range("K3").Value  = 430.9
$ws.Range("L3").Value  = 0.3631077778714081
$ws.Range("M3").Value  = -0
$ws.Range("N3").Value  = -0
$ws.Range("O3").Value  = -0
$ws.Range("P3").Value  = -0
$ws.Range("Q3").Value  = -0
$ws.Range("R3").Value  = -0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value  = 2074.2
$ws.Range("V3").Value  = 0.3588705491539499
$ws.Range("W3").Value  = 0.1003329685426223
$ws.Range("X3").Value  = 0.06367954879277794
$ws.Range("Y3").Value  = 0.03665341974984436
$ws.Range("Z3").Value  = 0.09230063234527763
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.03898376677617128
$ws.Range("AC3").Value = -0.03898376677617128
$ws.Range("AD3").Value = 8123.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 8123.8
$ws.Range("AG3").Value = 6049.6
$ws.Range("AH3").Value = 0.5842947150378319
$ws.Range("AI3").Value = 0.6199528384679369
$ws.Range("AJ3").Value = 0.511403790555734
$ws.Range("AK3").Value = 0.5484827329845779
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# --- Row 4: now MONETA Money Bank, a.s. (SEP:MONET) (was Komercni banka) --
$ws.Range("B4").Value  = "MONETA Money Bank, a.s. (SEP:MONET)"
$ws.Range("D4").Value  = -0.0931
$ws.Range("E4").Value  = -0.09080000000000001
$ws.Range("F4").ClearContents()
$ws.Range("K4").Value  = 124.4
$ws.Range("L4").Value  = 0.3766273085074175
$ws.Range("M4").Value  = 73.2
$ws.Range("N4").Value  = 0.04522426788582726
$ws.Range("O4").Value  = 0.5884244372990354
$ws.Range("P4").Value  = 73.2
$ws.Range("Q4").Value  = 0.04522426788582726
$ws.Range("R4").Value  = 0.5884244372990354
$ws.Range("T4").Value  = 0
$ws.Range("U4").Value  = 458.3
$ws.Range("V4").Value  = 0.2831459285802546
$ws.Range("W4").Value  = 0.1169502679326878
$ws.Range("X4").Value  = 0.05338436768340436
$ws.Range("Y4").Value  = 0.06356590024928342
$ws.Range("Z4").Value  = 0.1757101819342483
$ws.Range("AB4").Value = 0.0402297860858255
$ws.Range("AC4").Value = -0.0402297860858255
$ws.Range("AD4").Value = 1466
$ws.Range("AF4").Value = 1466
$ws.Range("AG4").Value = 1007.7
$ws.Range("AH4").Value = 0.4752642157816249
$ws.Range("AI4").Value = 0.5612987211884524
$ws.Range("AJ4").Value = 0.3836956935612839
$ws.Range("AK4").Value = 0.4679359182725795
